$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.347193598747253
$ws.Range("B1").Value = 1.725699663162231
$ws.Range("C1").Value = 2.482172250747681
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.208254337310791
